$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price record was added to the daily log: insert a fresh row right
# after the current row 96 (before the old row 97), which pushes every
# following row down by one (old row 97 -> new row 98, ..., old row 154 ->
# new row 155).
$ws.Rows.Item(97).Insert()

# Populate the newly inserted row 97 with the new record.
$ws.Cells.Item(97, 1).Value = 2
$ws.Cells.Item(97, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(97, 3).Value = "Coquimbo"
$ws.Cells.Item(97, 4).Value = 45119
$ws.Cells.Item(97, 5).Value = 4
$ws.Cells.Item(97, 6).Value = "Fruta"
$ws.Cells.Item(97, 7).Value = 100107
$ws.Cells.Item(97, 8).Value = "Otros"
$ws.Cells.Item(97, 9).Value = 100107011
$ws.Cells.Item(97, 10).Value = "Tuna"
$ws.Cells.Item(97, 11).Value = "Sin especificar"
$ws.Cells.Item(97, 12).Value = "Primera"
$ws.Cells.Item(97, 13).Value = 40
$ws.Cells.Item(97, 14).Value = 17000
$ws.Cells.Item(97, 15).Value = 18000
$ws.Cells.Item(97, 16).Value = 17500
$ws.Cells.Item(97, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(97, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(97, 19).Value = 972
$ws.Cells.Item(97, 20).Value = 18
